$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, which shifts existing data (A:D) to (B:E)
$ws.Range("A1").EntireColumn.Insert()

# Set new column widths (closest achievable values given Excel's internal
# column-width pixel snapping; target widths are 54.552101 / 6.596372 / 9.593605)
$ws.Columns.Item(1).ColumnWidth = 53.666666666666664
$ws.Columns.Item(2).ColumnWidth = 5.833333333333334
$ws.Columns.Item(3).ColumnWidth = 8.833333333333332
$ws.Columns.Item(4).ColumnWidth = 8.833333333333332
$ws.Columns.Item(5).ColumnWidth = 8.833333333333332

# Insert a new row before row 1, which shifts the data down by one row
$ws.Range("A1").EntireRow.Insert()

# Header row
$ws.Range("B1").Value = "Valid"
$ws.Range("C1").Value = "T"
$ws.Range("D1").Value = "Z"
$ws.Range("E1").Value = "p-value"

# Row labels in column A
$ws.Range("A2").Value = "CyclomaticComplexity(CC) & NbOperators"
$ws.Range("A3").Value = "CyclomaticComplexity(CC) & EffortToImplement"
$ws.Range("A4").Value = "MaintainabilityIndex & MaintainabilityIndex"
$ws.Range("A5").Value = "NbUniqueOperands & NbUniqueOperands"
$ws.Range("A6").Value = "NbOperands & NbOperands"
$ws.Range("A7").Value = "NbOperands & EffortToImplement"
$ws.Range("A8").Value = "NbUniqueOperators & ProgramLength"
$ws.Range("A9").Value = "NbOperators & CyclomaticComplexity(CC)"
$ws.Range("A10").Value = "NbOperators & NbOperators"
$ws.Range("A11").Value = "NbOperators & EffortToImplement"
$ws.Range("A12").Value = "ProgramLength & NbUniqueOperators"
$ws.Range("A13").Value = "ProgramLength & ProgramLength"
$ws.Range("A14").Value = "VocabularySize & VocabularySize"
$ws.Range("A15").Value = "ProgramVolume & ProgramVolume"
$ws.Range("A16").Value = "DifficultyLevel & DifficultyLevel"
$ws.Range("A17").Value = "ProgramLevel & ProgramLevel"
$ws.Range("A18").Value = "EffortToImplement & CyclomaticComplexity(CC)"
$ws.Range("A19").Value = "EffortToImplement & NbOperands"
$ws.Range("A20").Value = "EffortToImplement & NbOperators"
$ws.Range("A21").Value = "EffortToImplement & EffortToImplement"
$ws.Range("A22").Value = "TimeToImplement & TimeToImplement"

